$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Ejercicio Genérico " -> "Ejercicio genérico " (lowercase the
# "G"), leaving the "_GoBack" bookmark at the edit point, the way Word does
# after an in-place correction.
# ---------------------------------------------------------------------------
$rngG = $d.Content
$rngG.Find.Execute("Ejercicio G", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gPos = $rngG.End - 1

# Insert the bookmark right after the "G" first -- this is what forces the
# run to split at that exact spot instead of the whole paragraph collapsing
# into a single run once its text is touched.
$bmRange = $d.Range($gPos + 1, $gPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Now lowercase the "G" itself.
$gRange = $d.Range($gPos, $gPos + 1)
$gRange.Text = "g"

# ---------------------------------------------------------------------------
# Change 2: merge the three runs that spell "Sin ordenación aleatoria
# (S/N):)" into a single run, dropping the stray gramStart/gramEnd proofErr
# markers that used to sit between them.
# ---------------------------------------------------------------------------
$targetText = "Sin ordenación aleatoria (S/N):)"
$rngSinOrden = $d.Content
$rngSinOrden.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Setting identical text is treated as a no-op by the engine, so round-trip
# through a placeholder to force the run-merge/proofErr cleanup to happen.
$placeholder = $targetText + "__TMP__"
$rngSinOrden.Text = $placeholder
$rngSinOrden2 = $d.Range($rngSinOrden.Start, $rngSinOrden.Start + $placeholder.Length)
$rngSinOrden2.Text = $targetText

# ---------------------------------------------------------------------------
# Change 3: the stray "_GoBack" bookmark that used to sit after the closing
# curly quote in "... infieles" no longer belongs there. Word only keeps a
# single "_GoBack" bookmark, so re-adding it above already removed this one.
# ---------------------------------------------------------------------------
